{"js": "const replacements = [\n  [\"60\u00f79=\", \"36\u00f76=\"],\n  [\"22\u00f79=\", \"56\u00f73=\"],\n  [\"46\u00f74=\", \"11\u00f76=\"],\n  [\"31\u00f76=\", \"64\u00f78=\"],\n  [\"66\u00f79=\", \"63\u00f73=\"],\n  [\"39\u00f79=\", \"64\u00f75=\"],\n  [\"69\u00f74=\", \"22\u00f78=\"],\n  [\"71\u00f73=\", \"68\u00f75=\"],\n  [\"61\u00f79=\", \"60\u00f74=\"],\n  [\"93\u00f78=\", \"88\u00f76=\"],\n  [\"59\u00f73=\", \"32\u00f74=\"],\n  [\"90\u00f72=\", \"97\u00f73=\"],\n  [\"54\u00f74=\", \"18\u00f75=\"],\n  [\"17\u00f77=\", \"72\u00f72=\"],\n  [\"11\u00f72=\", \"27\u00f76=\"],\n  [\"19\u00f75=\", \"45\u00f75=\"],\n  [\"76\u00f79=\", \"48\u00f79=\"],\n  [\"57\u00f73=\", \"37\u00f74=\"],\n  [\"52\u00f74=\", \"15\u00f79=\"],\n  [\"45\u00f78=\", \"51\u00f74=\"],\n  [\"84\u00f72=\", \"25\u00f73=\"],\n  [\"91\u00f77=\", \"15\u00f79=\"],\n  [\"86\u00f75=\", \"44\u00f72=\"],\n  [\"76\u00f74=\", \"37\u00f74=\"],\n  [\"58\u00f75=\", \"80\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  // Only replace the first occurrence to preserve 1:1 mapping with the diff.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"60\u00f79=\", \"36\u00f76=\"),\n    @(\"22\u00f79=\", \"56\u00f73=\"),\n    @(\"46\u00f74=\", \"11\u00f76=\"),\n    @(\"31\u00f76=\", \"64\u00f78=\"),\n    @(\"66\u00f79=\", \"63\u00f73=\"),\n    @(\"39\u00f79=\", \"64\u00f75=\"),\n    @(\"69\u00f74=\", \"22\u00f78=\"),\n    @(\"71\u00f73=\", \"68\u00f75=\"),\n    @(\"61\u00f79=\", \"60\u00f74=\"),\n    @(\"93\u00f78=\", \"88\u00f76=\"),\n    @(\"59\u00f73=\", \"32\u00f74=\"),\n    @(\"90\u00f72=\", \"97\u00f73=\"),\n    @(\"54\u00f74=\", \"18\u00f75=\"),\n    @(\"17\u00f77=\", \"72\u00f72=\"),\n    @(\"11\u00f72=\", \"27\u00f76=\"),\n    @(\"19\u00f75=\", \"45\u00f75=\"),\n    @(\"76\u00f79=\", \"48\u00f79=\"),\n    @(\"57\u00f73=\", \"37\u00f74=\"),\n    @(\"52\u00f74=\", \"15\u00f79=\"),\n    @(\"45\u00f78=\", \"51\u00f74=\"),\n    @(\"84\u00f72=\", \"25\u00f73=\"),\n    @(\"91\u00f77=\", \"15\u00f79=\"),\n    @(\"86\u00f75=\", \"44\u00f72=\"),\n    @(\"76\u00f74=\", \"37\u00f74=\"),\n    @(\"58\u00f75=\", \"80\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceOne = 1: only replace the first occurrence encountered,\n    # preserving the 1:1 in-document-order mapping described by the diff.\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
